$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H52").Value = 3875

$ws.Range("H132").Value = 7112.3335
$ws.Range("I132").Value = 4061.4333
$ws.Range("K132").Value = 12184.2999
$ws.Range("M132").Value = -9654.2999

$ws.Range("H135").Value = 1147.1923
$ws.Range("I135").Value = 531.9
$ws.Range("K135").Value = 4787.099999999999
$ws.Range("M135").Value = -2252.099999999999

$ws.Range("H137").Value = 4835.1113
$ws.Range("I137").Value = 13445.444
$ws.Range("K137").Value = 40336.33199999999
$ws.Range("M137").Value = -37786.33199999999

$ws.Range("H138").Value = 18870048
$ws.Range("I138").Value = 30304842
$ws.Range("K138").Value = 90914526
$ws.Range("M138").Value = -90909386

$ws.Range("H141").Value = 5878.4595
$ws.Range("J141").Value = 6518.769
$ws.Range("L141").Value = 19556.307
$ws.Range("N141").Value = -29916.307

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 4263.7793
$ws.Range("I32").Value = 4432.107
$ws.Range("K32").Value = 4432.107
$ws.Range("M32").Value = -4145.107

$ws.Range("H45").Value = 11893.866
$ws.Range("I45").Value = 27599.6
$ws.Range("J45").Value = 4041
$ws.Range("K45").Value = 27599.6
$ws.Range("L45").Value = 4041
$ws.Range("M45").Value = -27222.6
$ws.Range("N45").Value = -4795

$ws.Range("H48").Value = 129899
$ws.Range("J48").Value = 129899
$ws.Range("L48").Value = 129899
$ws.Range("N48").Value = -130667

$ws.Range("H132").Value = 1653.4546
$ws.Range("I132").Value = 1653.4546
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4960.3638
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2430.3638
$ws.Range("N132").ClearContents()

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H42").Value = 169899
$ws.Range("J42").Value = 169899
$ws.Range("L42").Value = 169899
$ws.Range("N42").Value = -170555

$ws.Range("H43").Value = 189899
$ws.Range("J43").Value = 189899
$ws.Range("L43").Value = 189899
$ws.Range("N43").Value = -190261

$ws.Range("H47").Value = 129899
$ws.Range("J47").Value = 129899
$ws.Range("L47").Value = 129899
$ws.Range("N47").Value = -130939

$ws.Range("H48").Value = 129899
$ws.Range("J48").Value = 129899
$ws.Range("L48").Value = 129899
$ws.Range("N48").Value = -130729

$ws.Range("H99").Value = 1491.4166
$ws.Range("I99").Value = 1484.7
$ws.Range("K99").Value = 1484.7
$ws.Range("M99").Value = 13.29999999999995

$ws.Range("H134").Value = 2561.4082
$ws.Range("I134").Value = 2344.5227
$ws.Range("K134").Value = 7033.5681
$ws.Range("M134").Value = -4498.5681

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H7").Value = 778.2857
$ws.Range("I7").Value = 891.6667
$ws.Range("K7").Value = 891.6667
$ws.Range("M7").Value = -778.6667

$ws.Range("H31").Value = 2515.2856
$ws.Range("I31").Value = 2342.2856
$ws.Range("J31").Value = 2601.7856
$ws.Range("K31").Value = 2342.2856
$ws.Range("L31").Value = 2601.7856
$ws.Range("M31").Value = -2047.2856
$ws.Range("N31").Value = -3191.7856

$ws.Range("H34").Value = 2515.2856
$ws.Range("I34").Value = 2342.2856
$ws.Range("J34").Value = 2601.7856
$ws.Range("K34").Value = 2342.2856
$ws.Range("L34").Value = 2601.7856
$ws.Range("M34").Value = -2140.2856
$ws.Range("N34").Value = -3005.7856

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H105").Value = 2477.32
$ws.Range("I105").Value = 1252.381
$ws.Range("K105").Value = 1252.381
$ws.Range("M105").Value = 494.6189999999999

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H2").Value = 111111220
$ws.Range("I2").Value = 155555650
$ws.Range("J2").Value = 129.5
$ws.Range("K2").Value = 155555650
$ws.Range("L2").Value = 129.5
$ws.Range("M2").Value = -155555537
$ws.Range("N2").Value = -355.5

$ws.Range("H43").Value = 15000
$ws.Range("I43").Value = 15000
$ws.Range("K43").Value = 15000
$ws.Range("M43").Value = -14849

$ws.Range("H52").Value = 69030
$ws.Range("J52").Value = 69030
$ws.Range("L52").Value = 69030
$ws.Range("N52").Value = -69548

$ws.Range("H122").Value = 5554.143
$ws.Range("I122").Value = 5067.4585
$ws.Range("K122").Value = 15202.3755
$ws.Range("M122").Value = -12752.3755

$ws.Range("H132").Value = 3036.5435
$ws.Range("I132").Value = 2913.3635
$ws.Range("K132").Value = 8740.0905
$ws.Range("M132").Value = -6210.0905

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 3337.3333
$ws.Range("I46").Value = 2147.6667
$ws.Range("J46").Value = 4527
$ws.Range("K46").Value = 2147.6667
$ws.Range("L46").Value = 4527
$ws.Range("M46").Value = -1959.6667
$ws.Range("N46").Value = -4903

$ws.Range("H50").Value = 47499.5
$ws.Range("J50").Value = 45000
$ws.Range("L50").Value = 45000
$ws.Range("N50").Value = -46274

$ws.Range("H55").Value = 449.29413
$ws.Range("I55").Value = 183.76923
$ws.Range("J55").Value = 1312.25
$ws.Range("K55").Value = 183.76923
$ws.Range("L55").Value = 1312.25
$ws.Range("M55").Value = -10.76922999999999
$ws.Range("N55").Value = -1658.25

$ws.Range("H61").Value = 4134.1665
$ws.Range("I61").Value = 3198.5
$ws.Range("J61").Value = 6005.5
$ws.Range("K61").Value = 3198.5
$ws.Range("L61").Value = 6005.5
$ws.Range("M61").Value = -2996.5
$ws.Range("N61").Value = -6409.5

$ws.Range("H113").Value = 4134.1665
$ws.Range("I113").Value = 3198.5
$ws.Range("J113").Value = 6005.5
$ws.Range("K113").Value = 3198.5
$ws.Range("L113").Value = 6005.5
$ws.Range("M113").Value = -1028.5
$ws.Range("N113").Value = -10345.5

$ws.Range("H136").Value = 2493.125
$ws.Range("I136").Value = 1661.2593
$ws.Range("K136").Value = 4983.7779
$ws.Range("M136").Value = -2433.7779

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H81").Value = 2994
$ws.Range("I81").Value = 2994
$ws.Range("K81").Value = 5988
$ws.Range("M81").Value = -4927

$ws.Range("H84").Value = 2994
$ws.Range("I84").Value = 2994
$ws.Range("K84").Value = 29940
$ws.Range("M84").Value = -24636

$ws.Range("H96").Value = 115165.11
$ws.Range("I96").Value = 146898.28
$ws.Range("J96").Value = 4099
$ws.Range("K96").Value = 146898.28
$ws.Range("L96").Value = 4099
$ws.Range("M96").Value = -145525.28
$ws.Range("N96").Value = -6845

$ws.Range("H132").Value = 3986.0588
$ws.Range("I132").Value = 3430.1785
$ws.Range("K132").Value = 10290.5355
$ws.Range("M132").Value = -7760.5355

$ws.Range("H136").Value = 2713.0952
$ws.Range("I136").Value = 2813.75
$ws.Range("K136").Value = 8441.25
$ws.Range("M136").Value = -5891.25

Write-Output "Applied all changes"